$d = $word.ActiveDocument

$d.Content.Find.Execute("Your experience with the ParentText programme is vital to our study. We would love to hear about your experience of it, and what perhaps didn’t work for you. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. ", $true, $true, $false, $false, $false, $true, 1, $false, "Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva malunga namava akho ngayo, kwaye yintoni mhlawumbi engakhange ikusebenzele wena. This interview is part of a study carried out by researchers from the Universities of Cape Town in South Africa and the University of Oxford in the United Kingdom. ", 2) | Out-Null
$d.Content.Find.Execute("Before you decide if you’d like to be interviewed, it’s important for you to know why we’re doing this research and what participating in it would involve. All the information you might need is explained below, but if you have any questions about your participation or our study, please email the study team at swift@globalparenting.org or message us on WhatsApp at +27 XX XXX XXXX. We’re here to help you!", $true, $true, $false, $false, $false, $true, 1, $false, "Before you decide if you’d like to be interviewed, it’s important for you to know why we’re doing this research and what participating in it would involve. All the information you might need is explained below, but if you have any questions about your participation or our study, please email the study team at swift@globalparenting.org or message us on WhatsApp at +27 XX XXX XXXX. Silapha ukuzokunceda wena!", 2) | Out-Null
$d.Content.Find.Execute("What will my interview look like and what is expected of me?", $true, $true, $false, $false, $false, $true, 1, $false, "Udliwano-ndlebe lwam luzakujongeka njani kwaye kulindeleke ntoni kum?", 2) | Out-Null
$d.Content.Find.Execute("We would like to have a telephonic conversation with you which will last a maximum of 15 minutes. One of our researchers will call you to speak to you at a time that is convenient for you. There are no right or wrong answers, we just want to hear your experience and opinion of the chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. If while you are being interviewed, you are interrupted, please ask the researcher to pause until you feel safe to continue talking.", $true, $true, $false, $false, $false, $true, 1, $false, "We would like to have a telephonic conversation with you which will last a maximum of 15 minutes. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. Ukuba ngelixa wenziwa udliwano-ndlebe, uye waphazamiseka, ndicela ucele umphandi ukuba ame ude uzive ukhuselekile ukuqhubeka nokuthetha.", 2) | Out-Null
$d.Content.Find.Execute(" Why have I been invited to the interview?", $true, $true, $false, $false, $false, $true, 1, $false, " Kutheni ndimenyiwe nje kudliwano-ndlebe?", 2) | Out-Null
$d.Content.Find.Execute("We would like to know more about your experience with the chatbot. We see that you did not complete the programme and would like to know what we could have done differently to improve your experience, and improve the experience of other parents like you in the future. ", $true, $true, $false, $false, $false, $true, 1, $false, "We would like to know more about your experience with the chatbot. Siyabona ukuba akhange uyigqibe inkqubo kwaye singathanda ukwazi yintoni engesiyenze ngokwahlukileyo ukuphucula amava akho, kwaye siphucule namava abanye abazali abafana nawe kwixesha elizayo. ", 2) | Out-Null
$d.Content.Find.Execute("No, it's up to you if you want to be interviewed or not. If you don't want to be interviewed, there will be no implications to you or your family. If you do choose to be interviewed but don't want to answer some questions, you can skip these or stop at any time by telling your interviewer. ", $true, $true, $false, $false, $false, $true, 1, $false, "Hayi, kuxhomekeke kuwe ukuba uyafuna ukwenziwa udliwano-ndlebe okanye hayi. Ukuba awukufuni ukwenziwa udliwano-ndlebe, akuzokubakho ziphumo kuwe okanye kusapho lwakho. Ukuba uyakhetha ukwenziwa udliwano-ndlebe kodwa awufuni ukuphendula eminye imibuzo, ungayitsiba okanye ume nanini ngokuxelela umntu okwenza udliwano-ndlebe. ", 2) | Out-Null
$d.Content.Find.Execute(" What happens with my information?", $true, $true, $false, $false, $false, $true, 1, $false, " Kwenzeka ntoni ngolwazi lwam?", 2) | Out-Null
$d.Content.Find.Execute("To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview.", $true, $true, $false, $false, $false, $true, 1, $false, "Ukukhusela iinkcukacha zakho (kuquka igama lakho lokwenene, inkcukacha zoqhagamshelwano, kunye naluphi na olunye ulwazi olungakuchaza wena), sizakunika inombolo yokuthatha inxaxheba, kwaye ungazikhethela igama ofuna sikubize ngalo ngexesha lodliwano-ndlebe.", 2) | Out-Null
$d.Content.Find.Execute("With your permission, we will record the interview to help us remember the discussion and later write down what was said. We will delete any personal information we collect from you at the end of the study and, after transcribing your interview, change any data which might lead to identification at the point of transcription. We may use an artificial intelligence (AI) software, Microsoft Transcriber, to transcribe the interviews at first, and then we will check/review these transcriptions. This AI-generated information will be processed and stored securely on password-protected University of Cape Town servers, and in accordance with POPIA. Only authorised members of the research team will be able to access it, and this data will be owned by the Global Parenting Initiative at the University of Cape Town.", $true, $true, $false, $false, $false, $true, 1, $false, "With your permission, we will record the interview to help us remember the discussion and later write down what was said. Siza kucima nayiphi na ingcaciso yobuqu esiyiqokelele kuwe ekupheleni kophononongo kwaye, emva kokubhala udliwano-ndlebe lwakho, sitshintshe nayiphi na idatha enokukhokelela ekukuchazeni kwindawo yokukhuphela. Sinokusebenzisa i-software ye-Artificial Intelligence (AI), iMicrosoft Transcriber, ukukhuphela udliwano-ndlebe ekuqaleni, emva koko siya kujonga/sijongisise oku kukhutshelweyo. Olu lwazi luveliswe yi-AI luya kuqwalaselwa kwaye lugcinwe ngokukhuselekileyo kwiiseva zeDyunivesithi yaseKapa ezikhuselwe ngokuyimfihlo, kwaye ngokungqinelana nePOPIA. Ngamalungu eqela lophando kuphela agunyazisiweyo aya kukwazi ukufikelela kuyo, kwaye le datha iya kuba yeye Global Parenting Initiative kwiDyunivesithi yaseKapa.", 2) | Out-Null
$d.Content.Find.Execute(" Do I get anything for being interviewed? ", $true, $true, $false, $false, $false, $true, 1, $false, " Ingaba ikhona into endiyifumanayo ngokwenziwa oludliwano-ndlebe? ", 2) | Out-Null
$d.Content.Find.Execute("As a thank you for taking part in the discussion, we'll give you a R30 airtime voucher/data bundle. ", $true, $true, $false, $false, $false, $true, 1, $false, "Ukukubulela ngokuba yinxalenye kule ngxoxo, sizakunika R30 ye airtime voucher/data bundle. ", 2) | Out-Null
$d.Content.Find.Execute(" What happens to my information if I agree to be interviewed?", $true, $true, $false, $false, $false, $true, 1, $false, " Kwenzeka ntoni ngeenkcukacha zam ukuba ndiyavuma ukuba noludliwano-ndlebe?", 2) | Out-Null
$d.Content.Find.Execute("We only collect what’s needed for the study and store it securely. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. ", $true, $true, $false, $false, $false, $true, 1, $false, "Siqokelela kuphela oko sikudingayo koluphononongo kwaye sikugcina ngokukhuselekileyo. Your information, like your consent form and interview recording, and any information you provide via email or WhatsApp, will be kept safe on secure servers at the University of Cape Town. ", 2) | Out-Null
$d.Content.Find.Execute("Interview recordings will be deleted after we have written our notes. Any details that identify you will be kept separate and only authorised staff can access them. All data will be kept for five years after the study, but personal information will be deleted when the study ends. ", $true, $true, $false, $false, $false, $true, 1, $false, "Interview recordings will be deleted after we have written our notes. Nayiphi na inkcukacha echaza wena izakugcinwa bucala kwaye ngabasebenzi abagunyazisiweyo kuphela abanokufikelela kuzo. Yonke idatha iya kugcinwa iminyaka emihlanu emva koluphononongo, kodwa inkcukacha zomntu ziya kususwa xa isifundo siphelile. ", 2) | Out-Null
$d.Content.Find.Execute("Ethics committees and monitors may check the information. Your information will stay private unless the law says otherwise. After the study, we may share the information with other researchers but without your details. You have the right to see, correct, or ask us to delete your personal information.", $true, $true, $false, $false, $false, $true, 1, $false, "Iikomiti zokuziphatha kunye nabahloli banokulijonga ulwazi. Inkcukacha zakho ziyakuhlala ziyimfihlo ngaphandle kokuba umthetho uthetha enye into. Emva koluphononongo, singabelana ngolwazi nabanye abaphandi kodwa ngaphandle kweenkcukacha zakho. Unelungelo lokubona, ulungise, okanye ucele ukuba kucinywe ulwazi lwakho.", 2) | Out-Null
$d.Content.Find.Execute("You have the right to request access to your data, to correct any mistakes in your data, and to request us to delete it or transfer it somewhere else. Please email the study team before [*date to be determined] if you would like to do any of these.", $true, $true, $false, $false, $false, $true, 1, $false, "Unelungelo lokucela ukufikelela kwidatha yakho, ukulungisa naziphi na iimpazamo kwidatha yakho, kwaye usicele ukuba siyicime okanye siyidlulisele kwenye indawo. Nceda u-imeyilele iqela lophononongo phambi kwe [*umhla oza kumiselwa] ukuba ufuna ukwenza nayiphi na kwezi.", 2) | Out-Null
$d.Content.Find.Execute("Your participation and what you tell us will help us understand how to support families like yours. We plan to share the results in reports and at conferences so others can learn from this study too.", $true, $true, $false, $false, $false, $true, 1, $false, "Ukuthatha kwakho inxaxheba kunye nento osixelela yona izakusinceda siqondisise singazixhasa njani iintsapho ezifana nezakho. We plan to share the results in reports and at conferences so others can learn from this study too.", 2) | Out-Null
$d.Content.Find.Execute("Who are some of the study team members?", $true, $true, $false, $false, $false, $true, 1, $false, "Ngobani amanye amalungu eqela lophononongo?", 2) | Out-Null
